$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '316.56'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.18%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.07'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.46%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.183'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.19%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07612'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.84%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.321'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.03%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.652'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.13%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9318'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.60%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1252'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.26%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1825'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.75%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09081'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.92%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04128'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-3.49%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1055'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.53%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001271'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.98%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005978'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '5.17%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.358'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.20%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3362'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.46%'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.437'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '22.39%'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1362'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.86%'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2876'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '5.38%'
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04040'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.03%'
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001275'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.67%'
$ws.Range("B25").Value = 'HotbitToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004066'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.17%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001278'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.79%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02463'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '0.05%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05220'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-0.72%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007780'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.42%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1292'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.33%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007085'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4.23%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002167'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '17.83%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008212'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.50%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3437'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '11.10%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006693'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.98%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.84%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.3665'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '146.90%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004226'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '3.31%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002113'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.84%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002012'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.84%'
